$d = $word.ActiveDocument

# Update the date heading at the top of the document
$d.Content.Find.Execute("2026-01-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-06 Tuesday", 2) | Out-Null

# Update the arithmetic answers in the table, cell by cell (row, column)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "51-24=27"
$t.Cell(1,2).Range.Text = "34-15=19"
$t.Cell(1,3).Range.Text = "15+8=23"
$t.Cell(1,4).Range.Text = "85-29=56"
$t.Cell(1,5).Range.Text = "37+48=85"
$t.Cell(2,1).Range.Text = "8+26=34"
$t.Cell(2,2).Range.Text = "32-16=16"
$t.Cell(2,3).Range.Text = "90-11=79"
$t.Cell(2,4).Range.Text = "66+7=73"
$t.Cell(2,5).Range.Text = "35+29=64"
$t.Cell(3,1).Range.Text = "77+5=82"
$t.Cell(3,2).Range.Text = "31-18=13"
$t.Cell(3,3).Range.Text = "62-49=13"
$t.Cell(3,4).Range.Text = "78-49=29"
$t.Cell(3,5).Range.Text = "81-7=74"
$t.Cell(4,1).Range.Text = "62-45=17"
$t.Cell(4,2).Range.Text = "86-79=7"
$t.Cell(4,3).Range.Text = "76-17=59"
$t.Cell(4,4).Range.Text = "48+7=55"
$t.Cell(4,5).Range.Text = "49+9=58"
$t.Cell(5,1).Range.Text = "54-8=46"
$t.Cell(5,2).Range.Text = "91-29=62"
$t.Cell(5,3).Range.Text = "54-5=49"
$t.Cell(5,4).Range.Text = "7+67=74"
$t.Cell(5,5).Range.Text = "51-15=36"
$t.Cell(6,1).Range.Text = "81-29=52"
$t.Cell(6,2).Range.Text = "83-79=4"
$t.Cell(6,3).Range.Text = "81-14=67"
$t.Cell(6,4).Range.Text = "92-8=84"
$t.Cell(6,5).Range.Text = "93-39=54"
$t.Cell(7,1).Range.Text = "70-52=18"
$t.Cell(7,2).Range.Text = "36-18=18"
$t.Cell(7,3).Range.Text = "79+5=84"
$t.Cell(7,4).Range.Text = "51-8=43"
$t.Cell(7,5).Range.Text = "17+55=72"
$t.Cell(8,1).Range.Text = "5+79=84"
$t.Cell(8,2).Range.Text = "63+8=71"
$t.Cell(8,3).Range.Text = "60-44=16"
$t.Cell(8,4).Range.Text = "91-22=69"
$t.Cell(8,5).Range.Text = "36-19=17"
$t.Cell(9,1).Range.Text = "51-18=33"
$t.Cell(9,2).Range.Text = "7+39=46"
$t.Cell(9,3).Range.Text = "9+73=82"
$t.Cell(9,4).Range.Text = "26+35=61"
$t.Cell(9,5).Range.Text = "70-1=69"
$t.Cell(10,1).Range.Text = "37+56=93"
$t.Cell(10,2).Range.Text = "30-4=26"
$t.Cell(10,3).Range.Text = "95-87=8"
$t.Cell(10,4).Range.Text = "85-27=58"
$t.Cell(10,5).Range.Text = "61-14=47"
$t.Cell(11,1).Range.Text = "90-34=56"
$t.Cell(11,2).Range.Text = "61-12=49"
$t.Cell(11,3).Range.Text = "56+7=63"
$t.Cell(11,4).Range.Text = "81-53=28"
$t.Cell(11,5).Range.Text = "47+47=94"
$t.Cell(12,1).Range.Text = "75+17=92"
$t.Cell(12,2).Range.Text = "56+15=71"
$t.Cell(12,3).Range.Text = "12+49=61"
$t.Cell(12,4).Range.Text = "59+29=88"
$t.Cell(12,5).Range.Text = "23+19=42"
$t.Cell(13,1).Range.Text = "90-87=3"
$t.Cell(13,2).Range.Text = "73-4=69"
$t.Cell(13,3).Range.Text = "80-37=43"
$t.Cell(13,4).Range.Text = "7+77=84"
$t.Cell(13,5).Range.Text = "92-23=69"
$t.Cell(14,1).Range.Text = "47+19=66"
$t.Cell(14,2).Range.Text = "62-49=13"
$t.Cell(14,3).Range.Text = "78+7=85"
$t.Cell(14,4).Range.Text = "62-58=4"
$t.Cell(14,5).Range.Text = "66-19=47"
$t.Cell(15,1).Range.Text = "94-47=47"
$t.Cell(15,2).Range.Text = "67-39=28"
$t.Cell(15,3).Range.Text = "18+8=26"
$t.Cell(15,4).Range.Text = "35+56=91"
$t.Cell(15,5).Range.Text = "81-53=28"
$t.Cell(16,1).Range.Text = "7+36=43"
$t.Cell(16,2).Range.Text = "40-35=5"
$t.Cell(16,3).Range.Text = "20-15=5"
$t.Cell(16,4).Range.Text = "68-19=49"
$t.Cell(16,5).Range.Text = "46+38=84"
$t.Cell(17,1).Range.Text = "86-9=77"
$t.Cell(17,2).Range.Text = "58+37=95"
$t.Cell(17,3).Range.Text = "65-36=29"
$t.Cell(17,4).Range.Text = "24+9=33"
$t.Cell(17,5).Range.Text = "84-38=46"
$t.Cell(18,1).Range.Text = "77-68=9"
$t.Cell(18,2).Range.Text = "74-55=19"
$t.Cell(18,3).Range.Text = "9+5=14"
$t.Cell(18,4).Range.Text = "77-9=68"
$t.Cell(18,5).Range.Text = "73-55=18"
$t.Cell(19,1).Range.Text = "72-49=23"
$t.Cell(19,2).Range.Text = "26+48=74"
$t.Cell(19,3).Range.Text = "36+57=93"
$t.Cell(19,4).Range.Text = "30-19=11"
$t.Cell(19,5).Range.Text = "42-18=24"
$t.Cell(20,1).Range.Text = "81-64=17"
$t.Cell(20,2).Range.Text = "72+19=91"
$t.Cell(20,3).Range.Text = "44+7=51"
$t.Cell(20,4).Range.Text = "18+24=42"
$t.Cell(20,5).Range.Text = "82-67=15"
